# Update "想去人数" (interest count) values in column F for both the
# "展览" sheet and the "全部类型" sheet, reflecting newly generated data.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 802
$ws1.Range("F5").Value = 495
$ws1.Range("F6").Value = 1134
$ws1.Range("F9").Value = 114
$ws1.Range("F11").Value = 1151
$ws1.Range("F14").Value = 796
$ws1.Range("F16").Value = 186
$ws1.Range("F17").Value = 48
$ws1.Range("F20").Value = 194
$ws1.Range("F22").Value = 2368
$ws1.Range("F23").Value = 659
$ws1.Range("F24").Value = 70
$ws1.Range("F26").Value = 338
$ws1.Range("F27").Value = 2783
$ws1.Range("F32").Value = 103
$ws1.Range("F34").Value = 963
$ws1.Range("F36").Value = 337
$ws1.Range("F38").Value = 535

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 802
$ws4.Range("F6").Value = 495
$ws4.Range("F7").Value = 1134
$ws4.Range("F10").Value = 114
$ws4.Range("F12").Value = 1151
$ws4.Range("F14").Value = 796
$ws4.Range("F16").Value = 186
$ws4.Range("F19").Value = 48
$ws4.Range("F22").Value = 194
$ws4.Range("F24").Value = 2368
$ws4.Range("F25").Value = 659
$ws4.Range("F26").Value = 70
$ws4.Range("F29").Value = 2783
$ws4.Range("F39").Value = 103
$ws4.Range("F41").Value = 963
$ws4.Range("F44").Value = 337
$ws4.Range("F45").Value = 535
